# Agregando modulo de correos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2024-01", "Vj9W-c4Pm-ja0X-fC1C", 20157, 4998, 4526690, 0),
    @("2024-02", "Vj9W-c4Pm-ja0X-fC1C", 18768, 4737, 4253600, 0),
    @("2024-03", "Vj9W-c4Pm-ja0X-fC1C", 20518, 5045, 4588060, 0),
    @("2024-04", "Vj9W-c4Pm-ja0X-fC1C", 19861, 4966, 4472200, 0),
    @("2024-05", "Vj9W-c4Pm-ja0X-fC1C", 20352, 5093, 4559840, 0),
    @("2024-06", "Vj9W-c4Pm-ja0X-fC1C", 19705, 4913, 4441000, 0),
    @("2024-07", "Vj9W-c4Pm-ja0X-fC1C", 20387, 5029, 4565790, 0),
    @("2024-08", "Vj9W-c4Pm-ja0X-fC1C", 20358, 5030, 4560860, 0),
    @("2024-09", "Vj9W-c4Pm-ja0X-fC1C", 19907, 4838, 4481400, 0),
    @("2024-10", "Vj9W-c4Pm-ja0X-fC1C", 20440, 5137, 4574800, 0),
    @("2024-11", "Vj9W-c4Pm-ja0X-fC1C", 19651, 4793, 4430200, 0),
    @("2024-12", "Vj9W-c4Pm-ja0X-fC1C", 20281, 4991, 4547770, 0)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $rowIndex++
}
